# error-header-blank.xlsx: "services" sheet gets a new (whitespace-only)
# header cell in B1, and the active-sheet selection moves from A2 to B1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("services")

# New header cell - contains only spaces, used to exercise the
# "blank header" import-error test case.
$ws.Range("B1").Value = "       "

# Selection moves onto the newly added header cell.
$ws.Range("B1").Select()
